# Updates the PEBCOM interactive-map workbook:
# insert a brand new record as row 43 (pushing the existing rows 43-78 down
# to 44-79), matching the automatic daily map refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 43.
$ws.Rows.Item(43).Insert()

# Helper to write a value as literal text (preserving values such as
# "6076", "3" or "6/24/2025" as text instead of letting Excel reinterpret
# them as numbers/dates), without leaving a residual "Text" number format
# on the cell.
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value2 = "'" + $text
    $cell.Style = "Normal"
}

$r = 43

Set-TextCell $r 1 "6076"
Set-TextCell $r 2 "6/24/2025"
$ws.Cells.Item($r, 3).Value2 = "MATHEU 727"
Set-TextCell $r 4 "3"
Set-TextCell $r 5 "807763063"
$ws.Cells.Item($r, 6).Value2 = "PEBCOM"
$ws.Cells.Item($r, 7).Value2 = "Pendiente"
$ws.Cells.Item($r, 8).Value2 = "Colocar R400 para pedir a base traspaso de nodo propio y posterior a TLC"
$ws.Cells.Item($r, 9).Value2 = 1
$ws.Cells.Item($r, 10).Value2 = "Cambio"
$ws.Cells.Item($r, 11).Value2 = "Nodo TLC"
$ws.Cells.Item($r, 12).Value2 = "Pasante"
$ws.Cells.Item($r, 13).Value2 = -58.400169
$ws.Cells.Item($r, 14).Value2 = -34.617784
$ws.Cells.Item($r, 15).Value2 = "Almagro"
$ws.Cells.Item($r, 16).Value2 = "Capital Sur"

Write-Host "Inserted new row 43 (Caso 6076); sheet now spans to row $($ws.UsedRange.Rows.Count + $ws.UsedRange.Row - 1)."
